# Actualización desde MV -datos-
# Append the new daily "Dolar observado" rows (21-09-2021 .. 01-10-2021)
# to the bottom of the existing data table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(
    "21-09-2021",
    "22-09-2021",
    "23-09-2021",
    "24-09-2021",
    "27-09-2021",
    "28-09-2021",
    "29-09-2021",
    "30-09-2021",
    "01-10-2021"
)

$values = @(
    788.05,
    785.1,
    785.03,
    787.24,
    788.98,
    795.48,
    798.63,
    803.59,
    811.9
)

$startRow = 182
$endRow = $startRow + $dates.Count - 1

# Force column A of the new rows to be treated as plain text before typing
# the values in, so that day/month-ambiguous strings (like "01-10-2021")
# aren't auto-converted into date serial numbers by the smart text parser.
$dateRange = $ws.Range("A" + $startRow + ":A" + $endRow)
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $dates.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Drop the temporary text formatting again so the new cells end up with
# the same (default) style as all the other data rows in the sheet.
$dateRange.ClearFormats()
